# edit.ps1 -- apply the "literature_survey.docx" edits described by the
# commit's XML diff:
#   1. Collapse the per-letter "LITERATURE SURVEY" title runs into a
#      single run (keeps the same visible text).
#   2. Fix "forecastign" -> "forecasting" and add a comma before
#      "the two results".
#   3. Reword "A correlation with each attribute ..." paragraph.
#   4. Fix "Liner Regression" -> "Linear Regression".
#   5. Merge the "...reduce th" / "e dimensionality..." run split
#      (pure text merge, no visible change).
#   6. Reword "naive bayes" -> "naive Bayes," and drop "in order ".
#   7. Merge the "...L1 regular" / "ization..." run split (pure text
#      merge, no visible change).
#
# Helper: Word's Range.Find with wdReplaceAll (2) merges same-format
# runs it touches, but only actually rewrites XML when the replacement
# text differs from what's there already. For spots where the visible
# text does not change (pure run-merges) we nudge the range (drop the
# last character, then retype it) so the engine performs a genuine
# mutation and coalesces the backing runs.

$d = $word.ActiveDocument

function Replace-Text([string]$findText, [string]$replaceText) {
    $r = $d.Content
    $ok = $r.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        Write-Host "NOT FOUND: $findText"
    }
    return $ok
}

function Merge-Run([string]$fullText) {
    # Force a genuine (no-net-visible-change) mutation over $fullText so
    # the engine merges the underlying runs: shorten by the last
    # character, then restore it.
    $short = $fullText.Substring(0, $fullText.Length - 1)
    $lastChar = $fullText.Substring($fullText.Length - 1, 1)
    $r = $d.Content
    $ok = $r.Find.Execute($fullText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        Write-Host "NOT FOUND (merge): $fullText"
        return $false
    }
    $r.Text = $short
    $r2 = $d.Content
    $ok2 = $r2.Find.Execute($short, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok2) {
        Write-Host "NOT FOUND (merge restore): $short"
        return $false
    }
    $r2.Collapse(0)  # wdCollapseEnd
    $r2.InsertAfter($lastChar)
    return $true
}

# 1. Title: "LITERATURE SURVE" + "Y" -> collapse the 16 per-letter runs
#    that spell "LITERATURE SURVE" into one run (text is unchanged, so
#    use the merge helper; the trailing "Y" run is left as-is).
Merge-Run "LITERATURE SURVE"

# 2. "forecastign" typo + missing comma.
Replace-Text "price forecastign is done with externally generated residual value estimates and finally the two results" "price forecasting is done with externally generated residual value estimates and finally, the two results"

# 3. "A correlation with each attribute ..." paragraph reword.
Replace-Text "A correlation with each attribute to that of target attribute is found and linear regression curve with the target attribute is drawn. As a final step the total error and accuracy is measured." "A correlation between each attribute to that of the target attribute is found and a linear regression curve with the target attribute is drawn. As a final step, the total error and accuracy are measured."

# 4. "Liner Regression" -> "Linear Regression" (paper title).
Replace-Text "Car Price Prediction in the USA by using Liner Regression" "Car Price Prediction in the USA by using Linear Regression"

# 5. Pure run-merge: "...reduce th" / "e dimensionality..." (no visible
#    text change).
Merge-Run "to reduce the dimensionality"

# 6. "naive bayes" -> "naive Bayes," and drop "in order ".
Replace-Text "naïve bayes and decision trees have been used to make the predictions. The predictions are then evaluated and compared in order to find those which provide the best performances." "naïve Bayes, and decision trees have been used to make the predictions. The predictions are then evaluated and compared to find those which provide the best performances."

# 7. Pure run-merge: "...L1 regular" / "ization..." (no visible text
#    change). Keep inside the sentence, stop before the "models"
#    proofErr span that must stay untouched.
Merge-Run "10-fold cross-validation and L1 regularization. A general linear model, which"
